$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ODI Batting" sheet: remove the (empty) B3 cell entirely.
# ---------------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B3").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add a new worksheet "ODI Batting Extra" after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"

# Header row (text values).
$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

# Copy the header formatting (bold, border, centered) from an existing header cell.
$wsBatting.Range("A1").Copy()
$wsExtra.Range("A1:F1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsExtra.Application.CutCopyMode = $false

# Row 2 - match 4727: batting position/num4/num6/percent unknown (blank text cells).
$wsExtra.Range("A2").Value = "'4727"
$wsExtra.Range("A2").Style = "Normal"
$wsExtra.Range("B2").Value = "'"
$wsExtra.Range("B2").Style = "Normal"
$wsExtra.Range("C2").Value = "'"
$wsExtra.Range("C2").Style = "Normal"
$wsExtra.Range("D2").Value = "'"
$wsExtra.Range("D2").Style = "Normal"
$wsExtra.Range("E2").Value = "'"
$wsExtra.Range("E2").Style = "Normal"
$wsExtra.Range("F2").Value = "NO"

# Row 3 - match 4731: batting position is numeric (10).
$wsExtra.Range("A3").Value = "'4731"
$wsExtra.Range("A3").Style = "Normal"
$wsExtra.Range("B3").Value = 10
$wsExtra.Range("C3").Value = "'"
$wsExtra.Range("C3").Style = "Normal"
$wsExtra.Range("D3").Value = "'"
$wsExtra.Range("D3").Style = "Normal"
$wsExtra.Range("E3").Value = "'"
$wsExtra.Range("E3").Style = "Normal"
$wsExtra.Range("F3").Value = "NO"
